# Apply 2022-05-12 data update to "Fonds de solidarite volet 1" workbook.
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the rows whose underlying source figures changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => nombre_aides, nombre_entreprises, montant_total
$updates = @(
    @{ Row = 3;   C = 249324;  D = 47457;  E = 1036461686 },
    @{ Row = 6;   C = 20791;   D = 3261;   E = 360546521 },
    @{ Row = 7;   C = 7008;    D = 1200;   E = 290033265 },
    @{ Row = 63;  C = 14345;   D = 2812;   E = 36146424 },
    @{ Row = 64;  C = 5196;    D = 1053;   E = 20331329 },
    @{ Row = 70;  C = 15723;   D = 2768;   E = 24657676 },
    @{ Row = 74;  C = 939;     D = 150;    E = 4175586 },
    @{ Row = 79;  C = 116591;  D = 22734;  E = 447355948 },
    @{ Row = 91;  C = 151092;  D = 24834;  E = 482061222 },
    @{ Row = 92;  C = 408987;  D = 70904;  E = 1593473857 },
    @{ Row = 94;  C = 94135;   D = 13794;  E = 915375566 },
    @{ Row = 95;  C = 50704;   D = 6982;   E = 929563081 },
    @{ Row = 96;  C = 17226;   D = 2564;   E = 787987719 },
    @{ Row = 98;  C = 808;     D = 181;    E = 117608252 },
    @{ Row = 102; C = 107;     D = 23;     E = 19689236 },
    @{ Row = 105; C = 8168;    D = 1913;   E = 16867318 },
    @{ Row = 106; C = 18336;   D = 4979;   E = 41281761 },
    @{ Row = 107; C = 6386;    D = 1894;   E = 21935846 },
    @{ Row = 110; C = 394;     D = 71;     E = 16566525 },
    @{ Row = 111; C = 115;     D = 17;     E = 7718509 },
    @{ Row = 115; C = 11689;   D = 2248;   E = 32938514 },
    @{ Row = 142; C = 168971;  D = 35054;  E = 681755958 },
    @{ Row = 143; C = 64956;   D = 13059;  E = 373529758 },
    @{ Row = 145; C = 11831;   D = 1969;   E = 182701622 },
    @{ Row = 165; C = 83802;   D = 17112;  E = 354971515 },
    @{ Row = 166; C = 35929;   D = 7079;   E = 210575326 },
    @{ Row = 172; C = 22700;   D = 5292;   E = 44672707 },
    @{ Row = 175; C = 80778;   D = 14070;  E = 486144839 },
    @{ Row = 177; C = 14707;   D = 2256;   E = 251104307 },
    @{ Row = 178; C = 4922;    D = 815;    E = 212940437 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
